$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V content between row 63 and row 64 ---
# Capture row 63 values (F..V) before overwriting
$row63 = $ws.Range("F63:V63").Value2
$row64 = $ws.Range("F64:V64").Value2

$ws.Range("F63:V63").Value2 = $row64
$ws.Range("F64:V64").Value2 = $row63

# --- Append new row 83 ---
$ws.Cells.Item(83, 1).Value2 = 82
$ws.Cells.Item(83, 2).Value2 = "denmark"
$ws.Cells.Item(83, 3).Value2 = "superliga"
$ws.Cells.Item(83, 4).Value2 = "2023-2024"
$ws.Cells.Item(83, 5).Value2 = 45235.66666666666
$ws.Cells.Item(83, 6).Value2 = "Randers FC"
$ws.Cells.Item(83, 7).Value2 = 2
$ws.Cells.Item(83, 8).Value2 = "FC Copenhagen"
$ws.Cells.Item(83, 9).Value2 = 4
$ws.Cells.Item(83, 10).Value2 = 4.46
$ws.Cells.Item(83, 11).Value2 = "30/10/2023 19:13"
$ws.Cells.Item(83, 12).Value2 = 5.16
$ws.Cells.Item(83, 13).Value2 = "05/11/2023 15:52"
$ws.Cells.Item(83, 14).Value2 = 4
$ws.Cells.Item(83, 15).Value2 = "30/10/2023 19:13"
$ws.Cells.Item(83, 16).Value2 = 4.1
$ws.Cells.Item(83, 17).Value2 = "05/11/2023 15:52"
$ws.Cells.Item(83, 18).Value2 = 1.69
$ws.Cells.Item(83, 19).Value2 = "30/10/2023 19:13"
$ws.Cells.Item(83, 20).Value2 = 1.67
$ws.Cells.Item(83, 21).Value2 = "05/11/2023 15:49"
$ws.Cells.Item(83, 22).Value2 = "https://www.betexplorer.com/football/denmark/superliga/randers-fc-fc-copenhagen/dKa596LM/"

# --- Append new row 84 ---
$ws.Cells.Item(84, 1).Value2 = 83
$ws.Cells.Item(84, 2).Value2 = "denmark"
$ws.Cells.Item(84, 3).Value2 = "superliga"
$ws.Cells.Item(84, 4).Value2 = "2023-2024"
$ws.Cells.Item(84, 5).Value2 = 45235.75
$ws.Cells.Item(84, 6).Value2 = "Hvidovre IF"
$ws.Cells.Item(84, 7).Value2 = 1
$ws.Cells.Item(84, 8).Value2 = "Midtjylland"
$ws.Cells.Item(84, 9).Value2 = 4
$ws.Cells.Item(84, 10).Value2 = 5.55
$ws.Cells.Item(84, 11).Value2 = "29/10/2023 18:12"
$ws.Cells.Item(84, 12).Value2 = 8.26
$ws.Cells.Item(84, 13).Value2 = "05/11/2023 17:57"
$ws.Cells.Item(84, 14).Value2 = 4.27
$ws.Cells.Item(84, 15).Value2 = "29/10/2023 18:12"
$ws.Cells.Item(84, 16).Value2 = 5.09
$ws.Cells.Item(84, 17).Value2 = "05/11/2023 17:57"
$ws.Cells.Item(84, 18).Value2 = 1.58
$ws.Cells.Item(84, 19).Value2 = "29/10/2023 18:12"
$ws.Cells.Item(84, 20).Value2 = 1.39
$ws.Cells.Item(84, 21).Value2 = "05/11/2023 17:57"
$ws.Cells.Item(84, 22).Value2 = "https://www.betexplorer.com/football/denmark/superliga/hvidovre-if-midtjylland/M17hC8y4/"

# --- Apply styles matching the rest of the table ---
# Column A uses style index 1 (bold, bordered, centered) -- copy from existing row 82
$ws.Cells.Item(82, 1).Copy()
$ws.Cells.Item(83, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(82, 1).Copy()
$ws.Cells.Item(84, 1).PasteSpecial(-4122)

# Column E uses style index 2 (datetime number format) -- copy from existing row 82
$ws.Cells.Item(82, 5).Copy()
$ws.Cells.Item(83, 5).PasteSpecial(-4122)
$ws.Cells.Item(82, 5).Copy()
$ws.Cells.Item(84, 5).PasteSpecial(-4122)

$excel.CutCopyMode = $false
